# Add a "Generated TOC" column (G) to the learning-roadmap sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize existing columns E and F, new column G keeps the default width ---
$ws.Columns.Item(5).ColumnWidth = 28.5
$ws.Columns.Item(6).ColumnWidth = 34.8

# --- Header cell G1: copy the look of the neighboring header (F1) and set its text ---
$ws.Range("G1").Style = $ws.Range("F1").Style()
$ws.Range("G1").Value = "Generated TOC"

# --- Row 2 gets its own (non-shared) formula, matching how E2/F2 were authored ---
$ws.Range("G2").Formula = "=""<li><a href='""&B2&""/' target='_blank'>""&C2&"" Learning Resources</a></li>"""

# --- Rows 3-57 share one formula (fills down with relative references, like E3:E59/F3:F59) ---
$ws.Range("G3:G57").Formula = "=""<li><a href='""&B3&""/' target='_blank'>""&C3&"" Learning Resources</a></li>"""

# --- The trailing blank rows 58-59 no longer carry the old title/meta helper formulas ---
$ws.Range("E58:F59").ClearContents()

# --- Leave the selection on the newly generated TOC column, as in the final workbook ---
$ws.Range("G2:G57").Select() | Out-Null
